# Update "想去人数" (interest count) figures for several congan/exhibition
# events on both the "展览" and "全部类型" sheets, matching the new values
# published at the generated-site build 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12823
$ws1.Range("F6").Value = 320
$ws1.Range("F9").Value = 12818
$ws1.Range("F12").Value = 5224
$ws1.Range("F23").Value = 1154

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12823
$ws4.Range("F6").Value = 320
$ws4.Range("F10").Value = 12818
$ws4.Range("F13").Value = 5224
$ws4.Range("F25").Value = 1154
